$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

# New "Förändrad" (changed) date serial value: 46061 => 2026-02-08
$newDate = Get-Date -Year 2026 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = $newDate
    }
}
